$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at Q and R, shifting the old Q,R,S,T ("most_frequent_value",
# "memory_consumed_bytes", "pattern_count", "patterns") right to S,T,U,V.
$ws.Columns("Q:R").Insert()

# New header cells for the inserted columns.
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# New data values for every data row (2-10): default_count = 0,
# default_value = "<Unspecified>".
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 17).Value = 0
    $ws.Cells.Item($r, 18).Value = "<Unspecified>"
}

# The "most_frequent_value" column (now S) gets new literal values for several rows.
$ws.Range("S2").Value = ""
$ws.Range("S3").Value = "Private Hire Vehicle"
$ws.Range("S4").Value = "Licence Issued"
$ws.Range("S5").Value = "1/02/15 0:00"
$ws.Range("S6").Value = "31/01/16 0:00"
$ws.Range("S7").Value = "Skoda Octavia"
$ws.Range("S8").Value = ""
$ws.Range("S9").Value = "BLACK"
$ws.Range("S10").Value = "Amber Cars"
